# Fix missing/incorrect Cronbach's alpha values and the corresponding
# interpretation label (bug: some reliability values were stale/missing).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The alpha column holds text-like values (e.g. "0.791"), not numbers.
# Force the number format to Text first so Excel doesn't auto-convert
# the strings we assign into floating point numbers (which would also
# drop the trailing zero formatting, e.g. "0.780" -> 0.78).
$ws.Range("C2:C7").NumberFormat = "@"

$ws.Range("C2").Value = "0.780"
$ws.Range("C3").Value = "0.853"
$ws.Range("C4").Value = "0.887"
$ws.Range("C5").Value = "0.894"
$ws.Range("C6").Value = "0.652"
$ws.Range("C7").Value = "0.898"

$ws.Range("D7").Value = "Good"
